$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C3 with the new value (appended "modif 15h15" note)
$ws.Range("C3").Value = "Donnée C3 - modif 15h15"

# Mirror the author's last-edited cell becoming the active selection
$ws.Range("C3").Select()
